$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.015364243298753
$ws.Range("D2").Value = 1.021167648832802
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.013690346630008
$ws.Range("I2").Value = 1.026099249651183
$ws.Range("J2").Value = 1.02059018207206
$ws.Range("K2").Value = 1.024005777855235
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.016550741404369
$ws.Range("N2").Value = 1.022039537193933

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016346517302342
$ws.Range("D3").Value = 1.021849553635956
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.015316820857091
$ws.Range("I3").Value = 1.026215004597372
$ws.Range("J3").Value = 1.021207265213675
$ws.Range("K3").Value = 1.024494418912385
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.017979714334573
$ws.Range("N3").Value = 1.022657496664389

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.01698174162119
$ws.Range("D4").Value = 1.022290304323513
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.016368725254824
$ws.Range("I4").Value = 1.026288326665405
$ws.Range("J4").Value = 1.021605632708728
$ws.Range("K4").Value = 1.024809441117416
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.018903370668545
$ws.Range("N4").Value = 1.023056429886978

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017248701714699
$ws.Range("D5").Value = 1.022475479080598
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.016810825257778
$ws.Range("I5").Value = 1.026318773127447
$ws.Range("J5").Value = 1.021772885236777
$ws.Range("K5").Value = 1.02494159861179
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.019291446130043
$ws.Range("N5").Value = 1.023223919932801

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.017293520345345
$ws.Range("D6").Value = 1.022506563880959
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.016885048925344
$ws.Range("I6").Value = 1.026323863045662
$ws.Range("J6").Value = 1.021800954721672
$ws.Range("K6").Value = 1.024963772148413
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.019356592416429
$ws.Range("N6").Value = 1.023252029279584

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.016985309100398
$ws.Range("D7").Value = 1.022292779096991
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.016374633081875
$ws.Range("I7").Value = 1.026288734978456
$ws.Range("J7").Value = 1.021607868413436
$ws.Range("K7").Value = 1.024811208103873
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.018908557046742
$ws.Range("N7").Value = 1.023058668766643

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.015696285161701
$ws.Range("D8").Value = 1.021398202559479
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.014240135854391
$ws.Range("I8").Value = 1.026138696491979
$ws.Range("J8").Value = 1.020798920456006
$ws.Range("K8").Value = 1.024171156929384
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.017033877051458
$ws.Range("N8").Value = 1.022248572010327

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013421960570173
$ws.Range("D9").Value = 1.01981811846632
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.010474479731911
$ws.Range("I9").Value = 1.025862218915649
$ws.Range("J9").Value = 1.019366324389581
$ws.Range("K9").Value = 1.023034403253406
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.013722636993953
$ws.Range("N9").Value = 1.020813941493165

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01190372749183
$ws.Range("D10").Value = 1.018762229298096
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.007960641804139
$ws.Range("I10").Value = 1.025669775101214
$ws.Range("J10").Value = 1.018406424361843
$ws.Range("K10").Value = 1.02227056849797
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.0115095069206
$ws.Range("N10").Value = 1.019852678297287

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011245821567274
$ws.Range("D11").Value = 1.018304424818269
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.00687121680024
$ws.Range("I11").Value = 1.025584517903062
$ws.Range("J11").Value = 1.017989619303838
$ws.Range("K11").Value = 1.021938392802107
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.010549774613509
$ws.Range("N11").Value = 1.019435281328291

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011001369263941
$ws.Range("D12").Value = 1.018134285976328
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.006466409561581
$ws.Range("I12").Value = 1.025552559924964
$ws.Range("J12").Value = 1.01783462375059
$ws.Range("K12").Value = 1.021814792648913
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.010193064659216
$ws.Range("N12").Value = 1.019280065663579

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.011053808609866
$ws.Range("D13").Value = 1.018170785404346
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.00655324884642
$ws.Range("I13").Value = 1.025559428128024
$ws.Range("J13").Value = 1.017867878785494
$ws.Range("K13").Value = 1.021841315051707
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.010269590350368
$ws.Range("N13").Value = 1.019313367924447

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.011225616637922
$ws.Range("D14").Value = 1.0182903629226
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.006837758341595
$ws.Range("I14").Value = 1.02558188215579
$ws.Range("J14").Value = 1.017976810911104
$ws.Range("K14").Value = 1.021928180369712
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.01052029344249
$ws.Range("N14").Value = 1.01942245474617

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011331463019709
$ws.Range("D15").Value = 1.018364026665481
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.007013034389836
$ws.Range("I15").Value = 1.02559567843967
$ws.Range("J15").Value = 1.018043904290005
$ws.Range("K15").Value = 1.021981672415146
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.010674730165927
$ws.Range("N15").Value = 1.019489643405365

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.011947379859925
$ws.Range("D16").Value = 1.018792599676251
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.008032923311145
$ws.Range("I16").Value = 1.02567539272217
$ws.Range("J16").Value = 1.018434061772052
$ws.Range("K16").Value = 1.022292583735058
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.011573170379995
$ws.Range("N16").Value = 1.019880354955788

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012333593052569
$ws.Range("D17").Value = 1.019061272164724
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.008672421250976
$ws.Range("I17").Value = 1.025724879247246
$ws.Range("J17").Value = 1.018678485435702
$ws.Range("K17").Value = 1.022487226980453
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.012136349461215
$ws.Range("N17").Value = 1.020125125729077

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012558816436322
$ws.Range("D18").Value = 1.019217926639344
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.009045341719994
$ws.Range("I18").Value = 1.025753557902123
$ws.Range("J18").Value = 1.0188209415751
$ws.Range("K18").Value = 1.022600621128221
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.012464704588715
$ws.Range("N18").Value = 1.020267784172531

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012635603652881
$ws.Range("D19").Value = 1.019271331992748
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.009172483320445
$ws.Range("I19").Value = 1.025763305034828
$ws.Range("J19").Value = 1.018869496462648
$ws.Range("K19").Value = 1.022639262214249
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.012576642070843
$ws.Range("N19").Value = 1.020316408013589

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012292161031928
$ws.Range("D20").Value = 1.019032452112216
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.008603818322943
$ws.Range("I20").Value = 1.02571958905027
$ws.Range("J20").Value = 1.018652272686609
$ws.Range("K20").Value = 1.022466357875583
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.012075939971798
$ws.Range("N20").Value = 1.020098875754872

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011175025586508
$ws.Range("D21").Value = 1.018255152805768
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.006753981493816
$ws.Range("I21").Value = 1.025575278000565
$ws.Range("J21").Value = 1.017944737970733
$ws.Range("K21").Value = 1.021902606645791
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.010446473836521
$ws.Range("N21").Value = 1.019390336258545

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.010472193398289
$ws.Range("D22").Value = 1.017765913831739
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.0055900650119
$ws.Range("I22").Value = 1.02548286801974
$ws.Range("J22").Value = 1.017498866352309
$ws.Range("K22").Value = 1.021546907931925
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.009420671939626
$ws.Range("N22").Value = 1.01894383145128

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01084482069982
$ws.Range("D23").Value = 1.018025318002836
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.006207162671274
$ws.Range("I23").Value = 1.025532015167572
$ws.Range("J23").Value = 1.017735328025589
$ws.Range("K23").Value = 1.021735588800295
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.009964593886189
$ws.Range("N23").Value = 1.019180628927259

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012310882527225
$ws.Range("D24").Value = 1.01904547483133
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.008634817302627
$ws.Range("I24").Value = 1.025721980037199
$ws.Range("J24").Value = 1.018664117445482
$ws.Range("K24").Value = 1.022475788151842
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.012103236843493
$ws.Range("N24").Value = 1.020110737334661

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014010278481915
$ws.Range("D25").Value = 1.020227049017536
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.011448558876991
$ws.Range("I25").Value = 1.025935127405583
$ws.Range("J25").Value = 1.019737533872725
$ws.Range("K25").Value = 1.02332933762046
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.014579637246109
$ws.Range("N25").Value = 1.021185678136354
